$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 28: Tompkins health update - other data pending, quarantine no longer reported
$ws.Range("A28").Value = 26
$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 628
$ws.Range("P28").Value = 48
$ws.Range("Q28").Value = 515
$ws.Range("R28").Value = 1191
$ws.Range("S28").Value = 0
$ws.Range("T28").Value = 0
$ws.Range("U28").Value = 0
